$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 322 (pushes the existing rows 322-385 down to 325-388).
$ws.Rows("322:324").Insert()

# New row 322: Femacal de La Calera / Coquimbo / Frutilla / "Especial"
$ws.Cells.Item(322, 1).Value = 3
$ws.Cells.Item(322, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(322, 3).Value = "Coquimbo"
$ws.Cells.Item(322, 4).Value = 44889
$ws.Cells.Item(322, 5).Value = 5
$ws.Cells.Item(322, 6).Value = "Fruta"
$ws.Cells.Item(322, 7).Value = 100101
$ws.Cells.Item(322, 8).Value = "Berries"
$ws.Cells.Item(322, 9).Value = 100112025
$ws.Cells.Item(322, 10).Value = "Frutilla"
$ws.Cells.Item(322, 11).Value = "Sin especificar"
$ws.Cells.Item(322, 12).Value = "Especial"
$ws.Cells.Item(322, 13).Value = 150
$ws.Cells.Item(322, 14).Value = 8000
$ws.Cells.Item(322, 15).Value = 8500
$ws.Cells.Item(322, 16).Value = 8217
$ws.Cells.Item(322, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(322, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(322, 19).Value = 1174
$ws.Cells.Item(322, 20).Value = 7

# New row 323: same date, quality "Primera"
$ws.Cells.Item(323, 1).Value = 3
$ws.Cells.Item(323, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(323, 3).Value = "Coquimbo"
$ws.Cells.Item(323, 4).Value = 44889
$ws.Cells.Item(323, 5).Value = 5
$ws.Cells.Item(323, 6).Value = "Fruta"
$ws.Cells.Item(323, 7).Value = 100101
$ws.Cells.Item(323, 8).Value = "Berries"
$ws.Cells.Item(323, 9).Value = 100112025
$ws.Cells.Item(323, 10).Value = "Frutilla"
$ws.Cells.Item(323, 11).Value = "Sin especificar"
$ws.Cells.Item(323, 12).Value = "Primera"
$ws.Cells.Item(323, 13).Value = 97
$ws.Cells.Item(323, 14).Value = 6000
$ws.Cells.Item(323, 15).Value = 6000
$ws.Cells.Item(323, 16).Value = 6000
$ws.Cells.Item(323, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(323, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(323, 19).Value = 857
$ws.Cells.Item(323, 20).Value = 7

# New row 324: same date, quality "Segunda"
$ws.Cells.Item(324, 1).Value = 3
$ws.Cells.Item(324, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(324, 3).Value = "Coquimbo"
$ws.Cells.Item(324, 4).Value = 44889
$ws.Cells.Item(324, 5).Value = 5
$ws.Cells.Item(324, 6).Value = "Fruta"
$ws.Cells.Item(324, 7).Value = 100101
$ws.Cells.Item(324, 8).Value = "Berries"
$ws.Cells.Item(324, 9).Value = 100112025
$ws.Cells.Item(324, 10).Value = "Frutilla"
$ws.Cells.Item(324, 11).Value = "Sin especificar"
$ws.Cells.Item(324, 12).Value = "Segunda"
$ws.Cells.Item(324, 13).Value = 56
$ws.Cells.Item(324, 14).Value = 4000
$ws.Cells.Item(324, 15).Value = 4000
$ws.Cells.Item(324, 16).Value = 4000
$ws.Cells.Item(324, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(324, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(324, 19).Value = 571
$ws.Cells.Item(324, 20).Value = 7
